$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G5").Value = "AbeerRagheb@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G6").Value = "Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G18").Value = "ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G19").Value = "wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G20").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G21").Value = "hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G28").Value = "Sarah.Abdelmohsen@med.asu.edu.eg, dina.adel@med.asu.edu.eg, nourhan.osama@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg"
$ws.Range("G37").Value = "ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G38").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G39").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G40").Value = "hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G43").Value = "Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G44").Value = "Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G56").Value = "ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G57").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G58").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G59").Value = "asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G75").Value = "ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G76").Value = "wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G77").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G78").Value = "asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G94").Value = "ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G95").Value = "wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, eman.samir@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G96").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G97").Value = "asmaa.reda@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, heba@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G100").Value = "Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
$ws.Range("G113").Value = "ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G115").Value = "Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
